# version 3.01 description: 修改了sumsales_year下各个分类的列的前后顺序
#
# Swap the "分类名称" (category name) and "单品名称" (item name) columns:
# column A and column B exchange places, for the header row and every
# data row. All other columns (年份 / 销量(千克)) and values are untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row: A1 <-> B1
$ws.Cells.Item(1, 1).Value = "单品名称"
$ws.Cells.Item(1, 2).Value = "分类名称"

# Data rows 2..28: column A becomes the item name (was column B),
# column B becomes the category name (was column A, always "茄类").
$itemNames = @(
    "圆茄子(1)",
    "圆茄子(2)",
    "圆茄子(2)",
    "圆茄子(2)",
    "圆茄子(2)",
    "大龙茄子",
    "大龙茄子",
    "紫圆茄",
    "紫圆茄",
    "紫茄子(1)",
    "紫茄子(1)",
    "紫茄子(2)",
    "紫茄子(2)",
    "紫茄子(2)",
    "紫茄子(2)",
    "花茄子",
    "花茄子",
    "花茄子",
    "长线茄",
    "长线茄",
    "长线茄",
    "青茄子(1)",
    "青茄子(1)",
    "青茄子(1)",
    "青茄子(1)",
    "青茄子(2)",
    "青茄子(2)"
)

$categoryName = "茄类"

for ($i = 0; $i -lt $itemNames.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 1).Value = $itemNames[$i]
    $ws.Cells.Item($row, 2).Value = $categoryName
}
